$d = $word.ActiveDocument

# 1. Update the main title (Heading1) and the bolded "review" heading further down.
#    Both occurrences share identical text, so a single global replace covers them.
$d.Content.Find.Execute(
    "Play Dawn of Egypt Free - Our Review", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Dawn of Egypt Free - Exciting Slot Game Review", 2) | Out-Null

# 2. Update the "What we like" bullet list.
$d.Content.Find.Execute(
    "Stunning graphics and designs", $true, $false, $false, $false, $false,
    $true, 1, $false, "Stunning graphics", 2) | Out-Null

$d.Content.Find.Execute(
    "Egyptian themed story-line", $true, $false, $false, $false, $false,
    $true, 1, $false, "Exciting gameplay features", 2) | Out-Null

$d.Content.Find.Execute(
    "Free spin bonus and symbol upgrade feature", $true, $false, $false, $false, $false,
    $true, 1, $false, "Impressive winning symbols", 2) | Out-Null

$d.Content.Find.Execute(
    "Wide range of winning symbols", $true, $false, $false, $false, $false,
    $true, 1, $false, "Treasure trove for history lovers", 2) | Out-Null

# 3. Update the "What we don't like" bullet list, and add a brand new bullet
#    right after the existing one, matching its paragraph (ListBullet) style.
$d.Content.Find.Execute(
    "Limited paylines may discourage players", $true, $false, $false, $false, $false,
    $true, 1, $false, "Limited number of paylines", 2) | Out-Null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd("`r", "`a") -eq "Limited number of paylines") {
        $para.Range.InsertParagraphAfter()
        $newPara = $d.Paragraphs.Item($i + 1)
        $newPara.Range.Text = "Lack of innovative bonus features"
        break
    }
}

# 4. Update the meta title/description paragraphs near the end of the document.
$d.Content.Find.Execute(
    "Discover the stunning graphics and Egyptian theme of Dawn of Egypt slot game. Play for free with our review and enjoy its free spin bonus and symbol upgrade feature.",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Dawn of Egypt and play this exciting slot game for free. Experience stunning graphics and impressive winning symbols.", 2) | Out-Null
